$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-as-numbers to be written verbatim (matches source data which stores
# prices/volumes as text, not numeric values) by pre-marking numeric-looking cells
# as Text format before writing, same as typing an apostrophe-prefixed value in Excel.
$ws.Range("D2").Value = "61.760.31"
$ws.Range("E2").Value = "  -2.18%  "
$ws.Range("D3").Value = "3.394.02"
$ws.Range("E3").Value = "  -2.23%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "405.58"
$ws.Range("E5").Value = "  -2.40%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.37"
$ws.Range("E6").Value = "  +8.37%  "
$ws.Range("E7").Value = "  -1.15%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.669"
$ws.Range("E9").Value = "  -3.64%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.122"
$ws.Range("E10").Value = "  -7.89%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.43"
$ws.Range("E12").Value = "  -1.33%  "
$ws.Range("D13").Value = "3.917.10"
$ws.Range("E13").Value = "  -2.66%  "
$ws.Range("E14").Value = "  -2.67%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "19.83"
$ws.Range("E15").Value = "  -1.33%  "
$ws.Range("D16").Value = "3.399.27"
$ws.Range("E16").Value = "  -2.17%  "
$ws.Range("D17").Value = "61.669.09"
$ws.Range("E17").Value = "  -2.10%  "
$ws.Range("E18").Value = "  -1.69%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.00"
$ws.Range("E19").Value = "  +0.80%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000129"
$ws.Range("E20").Value = "  -8.31%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "3.19"
$ws.Range("E21").Value = "  -3.68%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "85.35"
$ws.Range("E22").Value = "  +4.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "314.73"
$ws.Range("E23").Value = "  -1.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.76"
$ws.Range("E24").Value = "  -0.98%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.13"
$ws.Range("E25").Value = "  -1.70%  "
$ws.Range("E26").Value = "  +11.10%  "
$ws.Range("E27").Value = "  -4.91%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.32"
$ws.Range("E28").Value = "  +5.68%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.67"
$ws.Range("E29").Value = "  -1.58%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.68"
$ws.Range("E30").Value = "  +4.24%  "
$ws.Range("E31").Value = "  -1.64%  "
$ws.Range("E32").Value = "  -0.56%  "
$ws.Range("B33").Value = "Cosmos"
$ws.Range("C33").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.36"
$ws.Range("E33").Value = "  -2.31%  "
$ws.Range("B34").Value = "Dai"
$ws.Range("C34").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  -0.56%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "41.56"
$ws.Range("E35").Value = "  -1.38%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0479"
$ws.Range("E36").Value = "  -2.72%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "51.78"
$ws.Range("E37").Value = "  -1.01%  "
$ws.Range("E38").Value = "  +0.15%  "
$ws.Range("E39").Value = "  -2.06%  "
$ws.Range("E40").Value = "  -3.06%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "139.30"
$ws.Range("E41").Value = "  +2.40%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.99"
$ws.Range("E42").Value = "  -0.75%  "
$ws.Range("E43").Value = "  -1.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.294"
$ws.Range("E44").Value = "  +4.35%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.99"
$ws.Range("E45").Value = "  +2.57%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.68"
$ws.Range("E46").Value = "  -1.40%  "
$ws.Range("E47").Value = "  -1.92%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "21.33"
$ws.Range("E48").Value = "  -3.11%  "
$ws.Range("D49").Value = "2.124.77"
$ws.Range("E49").Value = "  -2.69%  "
$ws.Range("E50").Value = "  -7.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.88"
$ws.Range("E51").Value = "  +0.65%  "
